$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 19-42 in column E ("UNIMAR HABIBLER(1)") hold the distances for the
# south-side lines; flip their sign to negative so the sheet can tell the
# "left"/"right" side lines apart (north side values, e.g. row 17-18, stay
# positive, as do rows whose column E cell is blank/"-").
for ($row = 19; $row -le 42; $row++) {
    $cell = $ws.Cells.Item($row, 5)  # Column E
    $val = $cell.Value2
    if (($val -ne $null) -and ($val -is [double]) -and ($val -gt 0)) {
        $cell.Value2 = -$val
    }
}
